$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Cells that need updated values while staying plain text (matching the
# original inline-string cells, i.e. no numeric conversion / no style change).
$changes = @{
    "D2" = "261.08"
    "E2" = "1.80%"
    "D3" = "27.36"
    "E3" = "1.66%"
    "D4" = "4.695"
    "E4" = "-2.29%"
    "E5" = "2.32%"
    "D6" = "6.675"
    "D7" = "0.8464"
    "E7" = "-0.48%"
    "D8" = "0.9281"
    "E8" = "0.32%"
    "D9" = "0.1404"
    "E9" = "1.76%"
    "D10" = "0.04818"
    "E10" = "15.26%"
    "D11" = "0.07100"
    "E11" = "1.29%"
    "D12" = "0.03083"
    "E12" = "1.04%"
    "D13" = "0.09063"
    "E13" = "-0.32%"
    "E14" = "0.42%"
    "D15" = "0.0006085"
    "E15" = "0.61%"
    "D16" = "0.006113"
    "E16" = "0.42%"
    "D17" = "3.450"
    "E17" = "-0.59%"
    "D18" = "3.147"
    "E18" = "-0.47%"
    "E21" = "0.19%"
    "D22" = "4.099"
    "E22" = "5.06%"
    "D23" = "0.04256"
    "E23" = "0.54%"
    "E24" = "0.24%"
    "D25" = "0.003795"
    "E25" = "4.85%"
    "D26" = "0.0001201"
    "E26" = "0.11%"
    "E27" = "3.45%"
    "D40" = "0.03873"
    "E40" = "2.41%"
    "E41" = "1.35%"
    "D42" = "0.004078"
    "E42" = "5.56%"
    "D43" = "0.01624"
    "E43" = "14.85%"
    "E44" = "-9.45%"
    "D45" = "0.00005148"
    "E45" = "-3.65%"
    "E46" = "0.08%"
    "D47" = "0.1371"
    "E47" = "-43.24%"
    "E48" = "23.80%"
    "E49" = "0.08%"
    "E50" = "0.08%"
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"          # force text interpretation so values
                                       # like "261.08" / "1.80%" are not
                                       # coerced into numbers/percentages
    $cell.Value = $changes[$addr]
    $cell.ClearFormats()               # drop the style index added by the
                                        # "@" format so the cell reverts to
                                        # having no explicit style, matching
                                        # the other (unchanged) data cells
}
